$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-9) get their per-record fields (Fecha, Volumen, Precio
# mínimo/máximo/promedio ponderado, Unidad de comercialización, Precio $/Kg)
# reshuffled among the rows. Columns A,B,C,E,F,G,H,I,J,K,L,R,T are identical
# across all rows and remain untouched.

$rows = @(
    @{ Row = 2;  D = 44208; M = 210; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 },
    @{ Row = 3;  D = 44351; M = 300; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 },
    @{ Row = 4;  D = 44162; M = 120; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 },
    @{ Row = 5;  D = 44397; M = 60;  N = 11000; O = 11000; P = 11000; Q = "`$/caja 14 kilos";           S = 786 },
    @{ Row = 6;  D = 44491; M = 180; N = 9000;  O = 9000;  P = 9000;  Q = "`$/caja 14 kilos empedrada"; S = 643 },
    @{ Row = 7;  D = 44400; M = 100; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos";           S = 714 },
    @{ Row = 8;  D = 44176; M = 250; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 },
    @{ Row = 9;  D = 44309; M = 300; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Volumen
    $ws.Cells.Item($row, 14).Value = $r.N   # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $r.O   # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $r.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($row, 19).Value = $r.S   # S: Precio $/Kg
}
